$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values must be forced to plain text (numeric-looking strings
# such as prices and percentages) so Excel does not auto-convert them to numbers.
$textUpdates = @{
    "D2" = "323.29"
    "E2" = "8.79%"
    "D3" = "49.54"
    "E3" = "18.67%"
    "D4" = "5.341"
    "E4" = "6.48%"
    "D5" = "0.08167"
    "E5" = "8.45%"
    "D6" = "4.599"
    "E6" = "5.11%"
    "D7" = "1.676"
    "E7" = "5.39%"
    "D8" = "1.207"
    "E8" = "31.07%"
    "D9" = "0.1346"
    "E9" = "13.09%"
    "D10" = "0.1970"
    "E10" = "7.64%"
    "D11" = "0.09747"
    "E11" = "8.71%"
    "D12" = "0.04495"
    "E12" = "9.84%"
    "E13" = "-0.17%"
    "D14" = "0.001326"
    "E14" = "3.54%"
    "D15" = "0.005981"
    "E15" = "3.73%"
    "D16" = "0.004272"
    "E16" = "9.10%"
    "D17" = "3.384"
    "E17" = "1.27%"
    "D18" = "2.437"
    "E18" = "1.51%"
    "D19" = "0.3394"
    "E19" = "1.93%"
    "D20" = "8.139"
    "E20" = "-3.10%"
    "D21" = "0.1418"
    "E21" = "2.67%"
    "D22" = "0.3050"
    "E22" = "-5.29%"
    "D23" = "0.04304"
    "E23" = "5.25%"
    "D24" = "0.001304"
    "E24" = "2.96%"
    "E25" = "9.58%"
    "D38" = "0.02763"
    "E38" = "14.72%"
    "D39" = "0.05592"
    "D40" = "0.006295"
    "E40" = "-0.14%"
    "D41" = "0.007682"
    "E41" = "-1.49%"
    "E42" = "9.40%"
    "D43" = "0.007677"
    "E43" = "3.76%"
    "D44" = "0.008097"
    "E44" = "3.81%"
    "D45" = "0.3522"
    "E45" = "18.74%"
    "D46" = "0.00006914"
    "E46" = "4.73%"
    "E47" = "-0.14%"
    "D48" = "0.06133"
    "E48" = "36.59%"
    "E49" = "-4.92%"
    "E50" = "-0.14%"
    "E51" = "-0.14%"
}

# Cells whose new values are plain (non-numeric-looking) text -- a simple
# assignment is sufficient since Excel will not reinterpret them as numbers.
$plainUpdates = @{
    "B16" = "HotbitToken"
    "C16" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "B18" = "BTSEToken"
    "C18" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B19" = "BitpandaEcosystemToken"
    "C19" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "B20" = "MCDex"
    "C20" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "B21" = "ProBitToken"
    "C21" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "B22" = "ZBToken"
    "C22" = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
    "B23" = "CoinExToken"
    "C23" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "B24" = "BitKan"
    "C24" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
}

foreach ($addr in $textUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Style = "Normal"
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$addr]
    $cell.Style = "Normal"
}

foreach ($addr in $plainUpdates.Keys) {
    $ws.Range($addr).Value = $plainUpdates[$addr]
}
